# Updated symbol list on Wed Feb  8 07:51:16 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) columns for
# the crypto symbols table on Sheet1 with newly scraped values. All of
# these columns hold plain text (e.g. "332.15", "1.11%") rather than
# numbers, so each new value is written with a leading apostrophe to force
# Excel to keep storing it as text instead of re-parsing it into a
# number/percentage and losing the original string formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'332.15"
$ws.Range("E2").Value  = "'1.11%"

$ws.Range("D3").Value  = "'45.83"
$ws.Range("E3").Value  = "'4.16%"

$ws.Range("D4").Value  = "'5.627"
$ws.Range("E4").Value  = "'2.38%"

$ws.Range("D5").Value  = "'0.08353"
$ws.Range("E5").Value  = "'4.33%"

$ws.Range("E6").Value  = "'2.78%"

$ws.Range("D7").Value  = "'0.9727"
$ws.Range("E7").Value  = "'2.45%"

$ws.Range("E8").Value  = "'-0.54%"

$ws.Range("D9").Value  = "'0.1165"
$ws.Range("E9").Value  = "'3.68%"

$ws.Range("D10").Value = "'0.1916"
$ws.Range("E10").Value = "'1.53%"

$ws.Range("E11").Value = "'-2.52%"

$ws.Range("D12").Value = "'0.09930"
$ws.Range("E12").Value = "'-1.56%"

$ws.Range("D13").Value = "'0.04704"
$ws.Range("E13").Value = "'-1.83%"

$ws.Range("D14").Value = "'0.1059"
$ws.Range("E14").Value = "'-0.30%"

$ws.Range("D15").Value = "'0.001288"
$ws.Range("E15").Value = "'1.11%"

$ws.Range("D16").Value = "'0.006034"
$ws.Range("E16").Value = "'0.73%"

$ws.Range("D17").Value = "'3.378"
$ws.Range("E17").Value = "'0.40%"

$ws.Range("D18").Value = "'4.452"
$ws.Range("E18").Value = "'1.81%"

$ws.Range("E19").Value = "'-3.10%"

$ws.Range("D20").Value = "'0.1393"
$ws.Range("E20").Value = "'-1.89%"

$ws.Range("D21").Value = "'0.2653"
$ws.Range("E21").Value = "'4.17%"

$ws.Range("D22").Value = "'0.04185"
$ws.Range("E22").Value = "'2.57%"

$ws.Range("D23").Value = "'0.001311"
$ws.Range("E23").Value = "'3.54%"

$ws.Range("D24").Value = "'0.004592"
$ws.Range("E24").Value = "'5.94%"

$ws.Range("D25").Value = "'0.0001301"
$ws.Range("E25").Value = "'8.37%"

$ws.Range("D26").Value = "'0.0003748"
$ws.Range("E26").Value = "'0.09%"

$ws.Range("D38").Value = "'0.02757"
$ws.Range("E38").Value = "'6.74%"

$ws.Range("D39").Value = "'0.05790"
$ws.Range("E39").Value = "'2.36%"

$ws.Range("D40").Value = "'0.007694"
$ws.Range("E40").Value = "'1.78%"

$ws.Range("D42").Value = "'0.007292"
$ws.Range("E42").Value = "'-1.52%"

$ws.Range("D43").Value = "'0.002012"
$ws.Range("E43").Value = "'-0.22%"

$ws.Range("D44").Value = "'0.008240"
$ws.Range("E44").Value = "'-4.55%"

$ws.Range("D45").Value = "'0.3403"

$ws.Range("D46").Value = "'0.00007287"
$ws.Range("E46").Value = "'2.50%"

$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.20%"

$ws.Range("D48").Value = "'0.0005813"
$ws.Range("E48").Value = "'0.02%"

$ws.Range("D49").Value = "'0.003508"
$ws.Range("E49").Value = "'-3.81%"

$ws.Range("D50").Value = "'0.003506"
$ws.Range("E50").Value = "'-0.71%"

$ws.Range("D51").Value = "'0.00002105"
$ws.Range("E51").Value = "'0.20%"
